# Insert a new data row at row 44 (shifts existing rows 44-80 down to 45-81)
# and populate it with the new record described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 44..80 down by one, creating a blank row 44.
$ws.Rows.Item(44).Insert(-4121)  # -4121 == xlShiftDown

# Populate the newly inserted row 44 with the new record's data.
$ws.Cells.Item(44, 1).Value = 11
$ws.Cells.Item(44, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(44, 3).Value = "Bíobío"
$ws.Cells.Item(44, 4).Value = 45049
$ws.Cells.Item(44, 5).Value = 8
$ws.Cells.Item(44, 6).Value = 100112031
$ws.Cells.Item(44, 7).Value = "Poroto verde"
$ws.Cells.Item(44, 8).Value = "Magnum"
$ws.Cells.Item(44, 9).Value = "Primera"
$ws.Cells.Item(44, 10).Value = 100
$ws.Cells.Item(44, 11).Value = 40000
$ws.Cells.Item(44, 12).Value = 42000
$ws.Cells.Item(44, 13).Value = 41000
$ws.Cells.Item(44, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(44, 15).Value = "Perú"
$ws.Cells.Item(44, 16).Value = 1640
$ws.Cells.Item(44, 17).Value = 25
$ws.Cells.Item(44, 18).Value = "Hortaliza"
